# 16-2-1.xlsx "Add files via upload" edit
# - Refresh the saved window geometry in the workbook view
# - Replace the data-reporter contact block (org / contact / email / phone / website)
#   with the new NSC Kyrgyz Republic details
# - Move the selection cursor to B8
# - Give B2 its own (Cyrillic-capable) font, matching the re-saved file's font table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect so we can write to the locked header/label cells,
# then restore protection exactly as it was.
$ws.Unprotect()

# --- Window geometry (workbook.xml bookViews) ---
$excel.Windows.Item(1).Left = 0
$excel.Windows.Item(1).Top = 0
$excel.Windows.Item(1).Width = 28800
$excel.Windows.Item(1).Height = 11835

# --- Data reporter block: new organization / contact / email / phone / website ---
$ws.Range("B6").Value = "National Statistical Committee of the Kyrgyz Republic (Department of Household Statistics)"
$ws.Range("B7").Value = "Kalymbetova Yryskan"
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com "
$ws.Range("B9").Value = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"

# --- Give B2 an explicit (Cyrillic-aware) font so it carries its own style ---
$ws.Range("B2").Font.Name = "Calibri"

# --- Selection cursor moves to B8 ---
$ws.Range("B8").Select()

# Restore sheet protection to its original state
$ws.Protect()
